$wb = $excel.ActiveWorkbook
$wsLib = $wb.Worksheets.Item("library_content")
$wsMes = $wb.Worksheets.Item("mesures")

# --- library_content sheet edits ---
# library_version: 1 -> 2
$wsLib.Range("B2").Value = 2

# library_description and framework_description: point to new "outline" text
$outlineText = "Swift Customer Security Controls Framework v2025 - outline"
$wsLib.Range("B6").Value = $outlineText
$wsLib.Range("B13").Value = $outlineText

# --- mesures sheet edits: remove "assessable" marker from non-leaf (category) rows ---
$wsMes.Range("A3").Clear()
$wsMes.Range("A11").Clear()
$wsMes.Range("A15").Clear()

# --- Restore view/selection state to match the edited file ---
$wsMes.Activate() | Out-Null
$wsMes.Range("C28").Select() | Out-Null

$wsLib.Activate() | Out-Null
$wsLib.Range("B3").Select() | Out-Null
